# Spain Primera Division RFEF - base update (28-06-2024 19:47)
#
# A handful of rows had their betting-odds records reassigned to a
# different match id. In the sheet layout, column A is a plain 0-based
# row counter (unchanged) while columns B..AD (id, Div, Date, teams,
# scores, odds, ...) carry the actual record. For each affected row we
# copy the full B..AD payload from another (source) row, forming closed
# permutation cycles. We snapshot every value first so the copy uses the
# original ("before") data even though some target rows are also used as
# a source for another target.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (source row's B..AD payload is copied into target row)
$rowMap = @{
    419 = 420
    420 = 419

    689 = 692
    691 = 689
    692 = 691

    709 = 711
    710 = 709
    711 = 710

    712 = 713
    713 = 712

    719 = 720
    720 = 719

    744 = 750
    745 = 748
    746 = 745
    747 = 746
    748 = 747
    749 = 744
    750 = 749

    753 = 755
    754 = 753
    755 = 760
    760 = 754
}

# Columns B (2) through AD (30) hold the record payload; column A (1) is
# left untouched.
$firstCol = 2
$lastCol = 30

# Snapshot every involved row's current (pre-edit) values before writing
# anything, since several rows are both a copy source and a copy target.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

foreach ($target in $rowMap.Keys) {
    $source = $rowMap[$target]
    $srcVals = $snapshot[$source]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($target, $c).Value = $srcVals[$c]
    }
}
